$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.071.20"
$ws.Range("E2").Value = "  +2.64%  "
$ws.Range("D3").Value = "1.679.84"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9949"
$ws.Range("E4").Value = "  -0.67%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.60"
$ws.Range("E5").Value = "  +7.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9961"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3660"
$ws.Range("E7").Value = "  +0.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.36"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3262"
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.152"
$ws.Range("E10").Value = "  +3.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07108"
$ws.Range("E11").Value = "  +2.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9966"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.110"
$ws.Range("E13").Value = "  +2.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.81"
$ws.Range("E14").Value = "  +3.61%  "
$ws.Range("D15").Value = "1.670.32"
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.642"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001057"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06619"
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9966"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "79.38"
$ws.Range("E20").Value = "  +3.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.01"
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.948"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.59"
$ws.Range("E23").Value = "  +3.33%  "
$ws.Range("D24").Value = "24.975.58"
$ws.Range("E24").Value = "  +2.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.455"
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.445"
$ws.Range("E26").Value = "  +4.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.87"
$ws.Range("E27").Value = "  +2.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.77"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").Value = "1.851.83"
$ws.Range("E29").Value = "  +1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.64"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.202"
$ws.Range("E31").Value = "  +4.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.069"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.804"
$ws.Range("E33").Value = "  +4.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08505"
$ws.Range("E34").Value = "  +2.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.641"
$ws.Range("E35").Value = "  -1.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.29"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.201"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02273"
$ws.Range("E38").Value = "  +2.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06094"
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.227"
$ws.Range("E40").Value = "  +2.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2091"
$ws.Range("E41").Value = "  +2.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.306"
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9964"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5964"
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.53"
$ws.Range("E45").Value = "  +6.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.838"
$ws.Range("E46").Value = "  +3.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5729"
$ws.Range("E47").Value = "  +2.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.92"
$ws.Range("E48").Value = "  +3.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.959"
$ws.Range("E49").Value = "  +1.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07023"
$ws.Range("E50").Value = "  +1.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.189"
$ws.Range("E51").Value = "  +3.86%  "
